$p = $ppt.ActivePresentation

# --- 1) Refresh the "datetimeFigureOut" date field from 19/11/2022 to 22/11/2022 ---
# This field lives on the slide master and on every slide layout's
# "Date Placeholder" shape. The only way to push new text into it through
# the exposed COM surface is TextRange.Text (PowerPoint itself collapses an
# auto-updating field to literal text the same way when scripted this way),
# so touch every placeholder currently showing the old date.

function Update-DateShape($shape) {
    if ($shape.HasTextFrame) {
        $tr = $shape.TextFrame.TextRange
        if ($tr.Text -eq "19/11/2022") {
            $tr.Text = "22/11/2022"
        }
    }
}

$master = $p.SlideMaster
foreach ($shp in $master.Shapes) {
    Update-DateShape $shp
}

foreach ($layout in $master.CustomLayouts) {
    foreach ($shp in $layout.Shapes) {
        Update-DateShape $shp
    }
}

# --- 2) Slide 3 ("User Interface"): rewrite the body placeholder text ---
$slide3 = $p.Slides.Item(3)
$body = $slide3.Shapes.Item(2)
$body.TextFrame.TextRange.Text = "Colors: Shades of Green and Red`r"
